# Insert a new weekly price-report row at row 178 (pushing existing rows
# 178..213 down to 179..214) and populate it with the new record's data.
# Everything below shifts down by one row, which Excel's row Insert does
# automatically (including copying the date-format style from the row
# above for column D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(178).Insert()

$ws.Cells.Item(178, 1).Value = 7
$ws.Cells.Item(178, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(178, 3).Value = "Ñuble"
$ws.Cells.Item(178, 4).Value = 44644
$ws.Cells.Item(178, 5).Value = 16
$ws.Cells.Item(178, 6).Value = 100112043
$ws.Cells.Item(178, 7).Value = "Pepino ensalada"
$ws.Cells.Item(178, 8).Value = "Sin especificar"
$ws.Cells.Item(178, 9).Value = "Primera"
$ws.Cells.Item(178, 10).Value = 120
$ws.Cells.Item(178, 11).Value = 18000
$ws.Cells.Item(178, 12).Value = 19000
$ws.Cells.Item(178, 13).Value = 18500
$ws.Cells.Item(178, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(178, 15).Value = "Región del Maule"
$ws.Cells.Item(178, 16).Value = 231
$ws.Cells.Item(178, 17).Value = 80
$ws.Cells.Item(178, 18).Value = "Hortaliza"
